# Se ajusta script de carga e insert
# Fills in missing "# Directivos Beneficiados" (E), "# Administrativos
# Beneficiados" (F), "¿Recibió Asistencia Técnica?" (I) and
# "¿Recibió Dotación?" (K) values for rows 2-26 of the staging sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ col letter = value }
$rowsData = @{
    2  = @{ E = 2; F = 0; I = "SI"; K = "NO" }
    3  = @{ E = 2; F = 0; I = "SI"; K = "NO" }
    4  = @{ E = 1; F = 0; I = "SI"; K = "NO" }
    5  = @{ E = 1; F = 0; I = "SI" }
    6  = @{ E = 3; F = 0; I = "SI" }
    7  = @{ E = 1; F = 0; I = "SI" }
    8  = @{          F = 0; I = "SI" }
    9  = @{ E = 1; F = 0; I = "SI" }
    10 = @{ E = 1; F = 0; I = "SI" }
    11 = @{ E = 1; F = 0; I = "SI" }
    12 = @{ E = 2; F = 0; I = "SI" }
    13 = @{ E = 2; F = 0; I = "SI" }
    14 = @{          F = 0; I = "SI" }
    15 = @{          F = 0; I = "SI" }
    16 = @{ E = 2; F = 0; I = "SI"; K = "NO" }
    17 = @{          F = 0; I = "SI" }
    18 = @{          F = 0; I = "SI" }
    19 = @{          F = 0; I = "SI" }
    20 = @{          F = 0; I = "SI" }
    21 = @{          F = 0; I = "SI" }
    22 = @{ E = 2; F = 0; I = "SI" }
    23 = @{ E = 2; F = 0; I = "SI" }
    24 = @{ E = 2; F = 0; I = "SI" }
    25 = @{          F = 0; I = "SI" }
    26 = @{          F = 0; I = "SI" }
}

foreach ($r in $rowsData.Keys) {
    $cols = $rowsData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
